$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "thin Oxea"
$ws.Range("A2").Value = "thick oxea"
